{"js": "// The author reworked the \"Firma\"/contribution column of the \"TAULA DE\n// CONTRIBUCIONS\" table: the plain \"Carles Maggi, Joan Maggi\" text in all\n// three data rows was rewritten to spell out the CM/JM initials used\n// elsewhere in the document: \"CM,JM (CM->Carles Maggi, JM->Joan Maggi)\".\nconst body = context.document.body;\n\nconst oldText = \"Carles Maggi, Joan Maggi\";\nconst newText = \"CM,JM (CM->Carles Maggi, JM->Joan Maggi)\";\n\nconst matches = body.search(oldText, { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < matches.items.length; i++) {\n  matches.items[i].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The author reworked the \"Firma\"/contribution column of the \"TAULA DE\n# CONTRIBUCIONS\" table: the plain \"Carles Maggi, Joan Maggi\" text in all\n# three data rows was rewritten to spell out the CM/JM initials used\n# elsewhere in the document: \"CM,JM (CM->Carles Maggi, JM->Joan Maggi)\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Carles Maggi, Joan Maggi\"\n$find.Replacement.Text = \"CM,JM (CM->Carles Maggi, JM->Joan Maggi)\"\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n"}
